# Apply updated crypto price / volume figures to the worksheet.
# Values mirror the inline-string cells already on the sheet (e.g. prices
# like '60.432.65' use '.' as a thousands separator and volume deltas like
# '  -1.69%  ' keep their padding spaces), so everything is written as literal
# text. Column D values get a leading apostrophe so Excel doesn't reinterpret
# numeric-looking strings (e.g. '136.14', '0.999') as actual numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    'D2' = '60.432.65'
    'D3' = '2.335.26'
    'D6' = '136.14'
    'D9' = '2.335.04'
    'D16' = '60.468.66'
    'D20' = '316.46'
    'D25' = '62.90'
    'D28' = '2.450.19'
    'D30' = '7.92'
    'D32' = '498.41'
    'D36' = '0.999'
    'D43' = '137.45'
    'D45' = '140.82'
}

$volumeUpdates = @{
    'E2' = '  -1.69%  '
    'E3' = '  -4.66%  '
    'E4' = '  -0.03%  '
    'E5' = '  -1.13%  '
    'E6' = '  -7.14%  '
    'E7' = '  +0.02%  '
    'E8' = '  -10.72%  '
    'E9' = '  -4.65%  '
    'E10' = '  -2.00%  '
    'E12' = '  -2.46%  '
    'E13' = '  -3.18%  '
    'E14' = '  -6.32%  '
    'E15' = '  -4.67%  '
    'E16' = '  -1.54%  '
    'E17' = '  -4.63%  '
    'E18' = '  -4.57%  '
    'E19' = '  -4.05%  '
    'E20' = '  -0.81%  '
    'E21' = '  -2.62%  '
    'E22' = '  -6.11%  '
    'E23' = '  -0.13%  '
    'E24' = '  -2.05%  '
    'E25' = '  -1.24%  '
    'E26' = '  +10.46%  '
    'E27' = '  +0.04%  '
    'E28' = '  -4.57%  '
    'E29' = '  -9.31%  '
    'E30' = '  -3.89%  '
    'E31' = '  -5.79%  '
    'E32' = '  -7.54%  '
    'E33' = '  -2.13%  '
    'E34' = '  -4.84%  '
    'E35' = '  -4.50%  '
    'E36' = '  -0.09%  '
    'E37' = '  -4.91%  '
    'E38' = '  -1.85%  '
    'E39' = '  -0.23%  '
    'E40' = '  -9.36%  '
    'E41' = '  +1.77%  '
    'E43' = '  -1.99%  '
    'E44' = '  -0.25%  '
    'E45' = '  -0.88%  '
    'E46' = '  -9.14%  '
    'E47' = '  -2.34%  '
    'E48' = '  -4.86%  '
    'E49' = '  -11.06%  '
    'E50' = '  -3.92%  '
    'E51' = '  -3.59%  '
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $priceUpdates[$addr]
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
